$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "255.42"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "4.04%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "28.19"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-5.08%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.199"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.43%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05861"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.95%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.696"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.65%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8701"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.38%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9599"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "12.51%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.96%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07160"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.81%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.03208"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.66%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09212"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.37%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.001551"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.10%"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "One"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0006059"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.93%"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005892"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.63%"
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.499"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.48%"
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.210"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.32%"
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.225"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.33%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.56%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.03463"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "3.78%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.528"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.99%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04173"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.82%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.76%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001225"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.01%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004560"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "9.21%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001200"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.05%"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001466"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "1.31%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03816"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.15%"
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1103"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.03%"
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.003838"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-33.51%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002354"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-2.58%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009732"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "6.02%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005413"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "2.35%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.11%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.09000"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "11.39%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002129"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-3.27%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.11%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002000"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.11%"
